$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sub_categories")

# Insert a new row above the current row 2 (shifts existing data down)
$ws.Rows.Item(2).Insert()

# Fill the new row 2 with "All" across columns A:I
$ws.Range("A2:I2").Value = "All"

# Update the active selection to match the new state
$ws.Range("I2").Select()
